# Rename the fund "Class/Series" labels from "Serie A"/"Serie B" to
# "Series A"/"Series B" as part of adding better validation for fund units.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Series A"
$ws.Range("B3").Value = "Series B"

# Move the active selection to B4 (matches the saved cursor position).
$ws.Range("B4").Select() | Out-Null
